$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.379.98'
$ws.Range('E2').Value = '  -1.86%  '
$ws.Range('D3').Value = '2.377.39'
$ws.Range('E3').Value = '  -1.19%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''502.33'
$ws.Range('E5').Value = '  -0.91%  '
$ws.Range('D6').Value = '''130.43'
$ws.Range('E6').Value = '  -2.21%  '
$ws.Range('D7').Value = '''0.999'
$ws.Range('E7').Value = '  +0.44%  '
$ws.Range('E8').Value = '  -1.80%  '
$ws.Range('D9').Value = '2.384.48'
$ws.Range('E9').Value = '  -2.38%  '
$ws.Range('D10').Value = '''0.0984'
$ws.Range('E10').Value = '  +0.50%  '
$ws.Range('E11').Value = '  +0.52%  '
$ws.Range('D12').Value = '''0.325'
$ws.Range('E12').Value = '  +0.84%  '
$ws.Range('E13').Value = '  +3.39%  '
$ws.Range('D14').Value = '2.802.88'
$ws.Range('E14').Value = '  -1.42%  '
$ws.Range('D15').Value = '56.351.25'
$ws.Range('E15').Value = '  -1.61%  '
$ws.Range('D16').Value = '''21.56'
$ws.Range('E16').Value = '  -1.65%  '
$ws.Range('E17').Value = '  -0.47%  '
$ws.Range('D18').Value = '2.406.21'
$ws.Range('E18').Value = '  +0.03%  '
$ws.Range('D19').Value = '''10.04'
$ws.Range('E19').Value = '  -2.45%  '
$ws.Range('E20').Value = '  -1.75%  '
$ws.Range('D21').Value = '''307.65'
$ws.Range('E21').Value = '  -1.98%  '
$ws.Range('D22').Value = '''6.27'
$ws.Range('E22').Value = '  -2.28%  '
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').Value = '''65.48'
$ws.Range('E24').Value = '  +0.52%  '
$ws.Range('E25').Value = '  +0.29%  '
$ws.Range('E26').Value = '  -3.29%  '
$ws.Range('E27').Value = '  -3.18%  '
$ws.Range('D28').Value = '''7.27'
$ws.Range('E28').Value = '  -4.17%  '
$ws.Range('D29').Value = '''171.73'
$ws.Range('E29').Value = '  -1.01%  '
$ws.Range('D30').Value = '0.0₃0716'
$ws.Range('E30').Value = '  -2.20%  '
$ws.Range('E31').Value = '  -2.47%  '
$ws.Range('B33').Value = 'Aptos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D33').Value = '''5.77'
$ws.Range('E33').Value = '  -6.63%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').Value = '''1.09'
$ws.Range('E34').Value = '  -3.47%  '
$ws.Range('E35').Value = '  +0.59%  '
$ws.Range('D36').Value = '''17.59'
$ws.Range('E36').Value = '  -1.95%  '
$ws.Range('D37').Value = '''1.17'
$ws.Range('E37').Value = '  -5.00%  '
$ws.Range('D38').Value = '''3.77'
$ws.Range('E38').Value = '  -1.36%  '
$ws.Range('E39').Value = '  -1.17%  '
$ws.Range('E40').Value = '  -2.15%  '
$ws.Range('D41').Value = '''1.38'
$ws.Range('E41').Value = '  -5.68%  '
$ws.Range('D42').Value = '''130.94'
$ws.Range('E42').Value = '  -2.88%  '
$ws.Range('D43').Value = '''3.36'
$ws.Range('E43').Value = '  -0.71%  '
$ws.Range('D44').Value = '''4.75'
$ws.Range('E44').Value = '  -4.89%  '
$ws.Range('D45').Value = '''0.565'
$ws.Range('E45').Value = '  -1.04%  '
$ws.Range('E46').Value = '  -0.88%  '
$ws.Range('D47').Value = '''241.31'
$ws.Range('E47').Value = '  -5.79%  '
$ws.Range('D48').Value = '''0.0482'
$ws.Range('E48').Value = '  -2.28%  '
$ws.Range('D49').Value = '''0.0208'
$ws.Range('E49').Value = '  -2.52%  '
$ws.Range('D50').Value = '''17.10'
$ws.Range('E50').Value = '  +0.53%  '
$ws.Range('E51').Value = '  -2.42%  '
